# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to reflect newly scraped counts.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1388
    $ws.Range("F3").Value = 2220
    $ws.Range("F4").Value = 348
    $ws.Range("F7").Value = 293
}
